# Adds an "Address" column (new column F) to the sheet, shifting the
# existing "District" column from F to G. Populates the new Address
# column for each data row (header + rows 3-50) with the school/taluk
# text; a handful of rows are intentionally left blank in column F
# where no address text exists in the source data (matching the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at F; existing F (District) shifts to G.
$ws.Columns("F:F").Insert()

# New column F values, keyed by row number (row 2 is the "Address" header).
$addresses = @{
    2 = "Address"
    3 = "Govt. H S ChintakiAurad"
    4 = "G H S BenchincholiHumnabad"
    5 = "Sri Veerbhadreshwar High School Chambol"
    6 = "Govt. Adarsha Vidayalaya Aurad"
    7 = "Govt. High School AlgoodBasavakalyan"
    8 = "Govt. High School(K) ManthalBasavakalyan"
    9 = "G H S Amlapur"
    10 = "Govt. HRPSBalnagar"
    11 = "Govt. Adarsh Vidyalaya (R M S A) Janwada"
    12 = "G H S HangaragaAurad"
    13 = "Smt. IndirabaiGurutappa ShetkarHigh School"
    14 = "Govt. High School AlgoodBasavakalyan"
    15 = "B P H S MadkattiBhalki"
    16 = "Govt. Urdu High School TalmadgiHumnabad"
    17 = "Govt. High School TornaAurad"
    18 = "G H S LadhaBhalki"
    19 = "Bhai Bansilal H S HalikhedHumanabad"
    20 = "Govt. High SchoolManthalBasavakalyan"
    23 = "Sri Babu Jagajeevan Ram High School KandgulAurad"
    24 = "G H S MurkiAurad"
    25 = "JIjamatGirls High School"
    26 = "Govt. High School Chitta(K)Basavakalyan"
    27 = "G H S Warwatti(B)Bhalki"
    28 = "Govt. High School EkambaAurad"
    29 = "Govt. Adarsh Vidyalaya RajolaBasavakalyan"
    30 = "R M H S HulsoorBasavakalyan"
    31 = "G H R P S Shamsheernagar"
    32 = "Ramswamy Periyar High SchoolKolar (K)"
    33 = "Basava High SchoolManhalli"
    34 = "G N P U C Basavakalyan"
    35 = "Govt. H S WalakhindiHumnabad"
    36 = "V B H S HulsoorBasavakalyan"
    37 = "Matoshri Ahilyabai HolkarHigh School"
    38 = "G H S YalladgundiBasavakalyan"
    39 = "Govt High School Khashampur"
    40 = "G H S HandikeraHumnabad"
    42 = "Govt. P U CollegeAurad(B)"
    43 = "Govt. HRPS YadalapurBasavakalyan"
    45 = "Govt. H S ChandikapurBasavakalyan"
    46 = "Humanabad"
    47 = "Sri. R. Govt. High School MuchalumBasavakalyan"
    49 = "Sri Basaveshwar Kannada H S Basavakalyan"
    50 = "Rural High School Sindol"
}

foreach ($row in $addresses.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $addresses[$row]
}
